# Rework the title block of the article into pandoc-style title/author
# paragraphs:
#   - paragraph 1 ("Day After Day - February 1935", style Heading1)
#     becomes style "Title"
#   - paragraph 2 ("By Dorothy Day", bold, no style) becomes style
#     "Authors" with text "Dorothy Day" (no bold)
#
# Each paragraph's text is rebuilt word-by-word (with separate runs for
# the inter-word spaces) to mirror how the upstream pandoc-based importer
# emits one run per token.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

function New-RunsXml($words) {
    $sb = ""
    foreach ($w in $words) {
        $escaped = $w.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $sb += "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    }
    return $sb
}

$titleWords = @("Day", " ", "After", " ", "Day", " ", "-", " ", "February", " ", "1935")
$authorWords = @("Dorothy", " ", "Day")

$titleRuns = New-RunsXml $titleWords
$authorRuns = New-RunsXml $authorWords

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$rng.InsertXML($xml)
